# AP150_TestData_ManageDistributionSets_21C.xlsx
# "Add files via upload" / "Anu - AP Files Uploaded"
#
# The uploaded workbook had the hard-coded Oracle Cloud login
# (URL / UserName / Password) stripped out of the Input_Value sheet
# before being committed, leaving the BR2:BT2 cells blank (formatting
# retained, values cleared).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Clear the credential values that used to live in BR2 (URL),
# BS2 (UserName) and BT2 (Password) - formatting/styles stay intact.
$ws.Range("BR2:BT2").ClearContents()

# Reflect the reviewer's on-screen selection at save time: scrolled
# further right and focused on the now-empty credential cells.
$ws.Range("BR2:BT2").Select()
